# Auto-generated Excel COM-interop script
# Applies numeric refresh updates to H..N columns across multiple sheets
# as captured by the scheduled runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 976.9091
$ws.Range("J6").Value = 891.5
$ws.Range("L6").Value = 2674.5
$ws.Range("N6").Value = -2898.5
$ws.Range("H17").Value = 3437.049
$ws.Range("J17").Value = 3472.6833
$ws.Range("L17").Value = 10418.0499
$ws.Range("N17").Value = -10754.0499
$ws.Range("H21").Value = 23749.834
$ws.Range("J21").Value = 29166.334
$ws.Range("L21").Value = 29166.334
$ws.Range("N21").Value = -30102.334
$ws.Range("H23").Value = 23749.834
$ws.Range("J23").Value = 29166.334
$ws.Range("L23").Value = 29166.334
$ws.Range("N23").Value = -29634.334
$ws.Range("H28").Value = 619.3333
$ws.Range("I28").Value = 459.5
$ws.Range("K28").Value = 459.5
$ws.Range("M28").Value = 25.5
$ws.Range("H33").Value = 311.9091
$ws.Range("J33").Value = 228
$ws.Range("L33").Value = 228
$ws.Range("N33").Value = -686
$ws.Range("H40").Value = 3250.25
$ws.Range("J40").Value = 3333.3333
$ws.Range("L40").Value = 3333.3333
$ws.Range("N40").Value = -3683.3333
$ws.Range("H76").Value = 33336334
$ws.Range("I76").Value = 50002500
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 50002500
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -50002185
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 33336334
$ws.Range("I79").Value = 50002500
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 50002500
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -50001408
$ws.Range("N79").Value = -6184
$ws.Range("H80").Value = 967.5
$ws.Range("J80").Value = 1216.8182
$ws.Range("L80").Value = 3650.4546
$ws.Range("N80").Value = -5646.4546
$ws.Range("H83").Value = 967.5
$ws.Range("J83").Value = 1216.8182
$ws.Range("L83").Value = 10951.3638
$ws.Range("N83").Value = -20935.3638
$ws.Range("H101").Value = 1863.2222
$ws.Range("I101").Value = 1972.25
$ws.Range("J101").Value = 1776
$ws.Range("K101").Value = 5916.75
$ws.Range("L101").Value = 5328
$ws.Range("M101").Value = -4294.75
$ws.Range("N101").Value = -8572
$ws.Range("H115").Value = 7636510
$ws.Range("I115").Value = 8181817.5
$ws.Range("J115").Value = 2204
$ws.Range("K115").Value = 24545452.5
$ws.Range("L115").Value = 6612
$ws.Range("M115").Value = -24543885.5
$ws.Range("N115").Value = -9746
$ws.Range("H137").Value = 5642.0435
$ws.Range("I137").Value = 2611.111
$ws.Range("J137").Value = 16553.4
$ws.Range("K137").Value = 7833.333
$ws.Range("L137").Value = 49660.2
$ws.Range("M137").Value = -5283.333
$ws.Range("N137").Value = -54760.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36901.64
$ws.Range("I32").Value = 36726.363
$ws.Range("J32").Value = 38829.668
$ws.Range("K32").Value = 36726.363
$ws.Range("L32").Value = 38829.668
$ws.Range("M32").Value = -36439.363
$ws.Range("N32").Value = -39403.668
$ws.Range("H61").Value = 12354470
$ws.Range("J61").Value = 12719.363
$ws.Range("L61").Value = 12719.363
$ws.Range("N61").Value = -13143.363
$ws.Range("H136").Value = 12354470
$ws.Range("J136").Value = 12719.363
$ws.Range("L136").Value = 38158.089
$ws.Range("N136").Value = -43258.089

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3092.1
$ws.Range("J4").Value = 4373
$ws.Range("L4").Value = 4373
$ws.Range("N4").Value = -4597
$ws.Range("H31").Value = 30306416
$ws.Range("I31").Value = 58825580
$ws.Range("J31").Value = 4803.4375
$ws.Range("K31").Value = 58825580
$ws.Range("L31").Value = 4803.4375
$ws.Range("M31").Value = -58825285
$ws.Range("N31").Value = -5393.4375
$ws.Range("H34").Value = 30306416
$ws.Range("I34").Value = 58825580
$ws.Range("J34").Value = 4803.4375
$ws.Range("K34").Value = 58825580
$ws.Range("L34").Value = 4803.4375
$ws.Range("M34").Value = -58825378
$ws.Range("N34").Value = -5207.4375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 207108.95
$ws.Range("I4").Value = 240447.23
$ws.Range("J4").Value = 7079.25
$ws.Range("K4").Value = 721341.6900000001
$ws.Range("L4").Value = 21237.75
$ws.Range("M4").Value = -721229.6900000001
$ws.Range("N4").Value = -21461.75
$ws.Range("H62").Value = 2028.4637
$ws.Range("I62").Value = 1235.8948
$ws.Range("K62").Value = 3707.6844
$ws.Range("M62").Value = -3021.6844
$ws.Range("H65").Value = 2028.4637
$ws.Range("I65").Value = 1235.8948
$ws.Range("K65").Value = 11123.0532
$ws.Range("M65").Value = -7691.0532
$ws.Range("H107").Value = 1726.5588
$ws.Range("I107").Value = 633.3333
$ws.Range("J107").Value = 2322.8635
$ws.Range("K107").Value = 1899.9999
$ws.Range("L107").Value = 6968.5905
$ws.Range("M107").Value = 20.00009999999997
$ws.Range("N107").Value = -10808.5905

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 64949.668
$ws.Range("I20").Value = 115000
$ws.Range("J20").Value = 39924.5
$ws.Range("K20").Value = 115000
$ws.Range("L20").Value = 39924.5
$ws.Range("M20").Value = -114755
$ws.Range("N20").Value = -40414.5
$ws.Range("H70").Value = 4378.4443
$ws.Range("I70").Value = 4425.75
$ws.Range("K70").Value = 4425.75
$ws.Range("M70").Value = -4155.75
$ws.Range("H73").Value = 4378.4443
$ws.Range("I73").Value = 4425.75
$ws.Range("K73").Value = 4425.75
$ws.Range("M73").Value = -3489.75
$ws.Range("H80").Value = 3384
$ws.Range("I80").Value = 2456.2856
$ws.Range("K80").Value = 2456.2856
$ws.Range("M80").Value = -1458.2856
$ws.Range("H83").Value = 3384
$ws.Range("I83").Value = 2456.2856
$ws.Range("K83").Value = 12281.428
$ws.Range("M83").Value = -7289.428
$ws.Range("H97").Value = 1086.2222
$ws.Range("I97").Value = 1059.5
$ws.Range("K97").Value = 1059.5
$ws.Range("M97").Value = -563.5
$ws.Range("H102").Value = 1190.3158
$ws.Range("I102").Value = 1204.8667
$ws.Range("K102").Value = 1204.8667
$ws.Range("M102").Value = 417.1333
$ws.Range("H113").Value = 1070.4
$ws.Range("I113").Value = 1070.4
$ws.Range("K113").Value = 1070.4
$ws.Range("M113").Value = 1099.6
$ws.Range("H118").Value = 29326.666
$ws.Range("J118").Value = 29326.666
$ws.Range("L118").Value = 29326.666
$ws.Range("N118").Value = -32640.666
$ws.Range("H132").Value = 4497.3105
$ws.Range("I132").Value = 2084.7778
$ws.Range("J132").Value = 8445.091
$ws.Range("K132").Value = 6254.3334
$ws.Range("L132").Value = 25335.273
$ws.Range("M132").Value = -3724.3334
$ws.Range("N132").Value = -30395.273

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 94015.63
$ws.Range("I7").Value = 169162.83
$ws.Range("K7").Value = 169162.83
$ws.Range("M7").Value = -169050.83
$ws.Range("H14").Value = 25000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 25000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 25000
# M14 is dropped entirely by the refresh (HQ profit no longer applicable)
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -25344
$ws.Range("H40").Value = 14324.071
$ws.Range("I40").Value = 11276.318
$ws.Range("K40").Value = 11276.318
$ws.Range("M40").Value = -11140.318
$ws.Range("H55").Value = 167.1875
$ws.Range("I55").Value = 159.33333
$ws.Range("J55").Value = 171.9
$ws.Range("K55").Value = 159.33333
$ws.Range("L55").Value = 171.9
$ws.Range("M55").Value = 13.66667000000001
$ws.Range("N55").Value = -517.9
$ws.Range("H61").Value = 3193.111
$ws.Range("I61").Value = 2139.8333
$ws.Range("J61").Value = 5299.6665
$ws.Range("K61").Value = 2139.8333
$ws.Range("L61").Value = 5299.6665
$ws.Range("M61").Value = -1937.8333
$ws.Range("N61").Value = -5703.6665
$ws.Range("H113").Value = 3193.111
$ws.Range("I113").Value = 2139.8333
$ws.Range("J113").Value = 5299.6665
$ws.Range("K113").Value = 2139.8333
$ws.Range("L113").Value = 5299.6665
$ws.Range("M113").Value = 30.16670000000022
$ws.Range("N113").Value = -9639.666499999999
$ws.Range("H122").Value = 3560.4
$ws.Range("I122").Value = 3451
$ws.Range("K122").Value = 10353
$ws.Range("M122").Value = -7903
$ws.Range("H126").Value = 94015.63
$ws.Range("I126").Value = 169162.83
$ws.Range("K126").Value = 507488.49
$ws.Range("M126").Value = -505018.49
$ws.Range("H136").Value = 1239605
$ws.Range("I136").Value = 1554611
$ws.Range("K136").Value = 4663833
$ws.Range("M136").Value = -4661283

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 1687.25
$ws.Range("I19").Value = 1250.3334
$ws.Range("J19").Value = 2998
$ws.Range("K19").Value = 1250.3334
$ws.Range("L19").Value = 2998
$ws.Range("M19").Value = -1076.3334
$ws.Range("N19").Value = -3346
$ws.Range("H29").Value = 356.25
$ws.Range("J29").Value = 691.6667
$ws.Range("L29").Value = 691.6667
$ws.Range("N29").Value = -1271.6667
$ws.Range("H101").Value = 25347.334
$ws.Range("J101").Value = 25347.334
$ws.Range("L101").Value = 25347.334
$ws.Range("N101").Value = -31837.334
